$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells remain text (matching original inlineStr type) rather than
# being auto-converted to numbers/dates by Excel when values look numeric.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.288.76"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.71%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.505.66"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.05%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.35"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.56"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.44%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.506.65"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.17%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.10"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.54%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.373"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.41%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.101.80"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000179"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.47%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.508.52"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.18"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -5.68%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.299.38"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.71%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.81"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.57"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "383.22"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.647.52"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.567"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.83%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.94"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.73"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +4.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.58"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.97%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.28%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.26"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.52%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.522.73"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.41%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "23.54"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.81%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.86%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.54"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.84"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.79%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "164.21"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -4.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0781"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.809"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.77%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "25.85"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.43%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "41.84"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.20%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.30%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.05%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.470.64"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.57%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.28%  "
